$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1. Title: capitalize "Funcionalidad" and "Diseño"
# -----------------------------------------------------------------
$d.Content.Find.Execute(
    "Informe de funcionalidad y diseño - Entrega Final Python",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Informe de Funcionalidad y Diseño - Entrega Final Python", 2) | Out-Null

# -----------------------------------------------------------------
# 2. Register a "Hyperlink" character style so the style reference
#    Word writes for the new hyperlink resolves cleanly.
# -----------------------------------------------------------------
try {
    $d.Styles.Add("Hyperlink", 2) | Out-Null
} catch {
}

# -----------------------------------------------------------------
# 3. Replace the old (empty) Subtitle paragraph with the new block
#    of paragraphs describing the GitHub repository + heading.
# -----------------------------------------------------------------
try {
    $bk = $d.Bookmarks("_q19sdmjuug63")
    $bk.Delete()
} catch {
}

$subtitlePara = $d.Paragraphs(2)
$subtitlePara.Range.Style = "Normal"
$subtitlePara.Range.Text = "`rEl repositorio GitHub solicitado se encuentra en el siguiente enlace:`r`r`rInforme de funcionalidad y diseño:`r"

# Paragraph 4 is now the (currently empty) paragraph meant to host the
# hyperlink to the GitHub repository.
$linkPara = $d.Paragraphs(4)
$linkRange = $linkPara.Range
$linkRange.Collapse(1)
$url = "https://github.com/guidojuant/Entrega_Final_Python_ITBA"
$hyperlink = $d.Hyperlinks.Add($linkRange, $url, $null, $null, $url)
$hyperlink.Range.Font.Color = 13391121
$hyperlink.Range.Font.Underline = 1

# Center the "Informe de funcionalidad y diseño:" heading paragraph.
$headingPara = $d.Paragraphs(6)
$headingPara.Alignment = 1

# -----------------------------------------------------------------
# 4. Fix "aplicacion" -> "aplicación" typo.
# -----------------------------------------------------------------
$d.Content.Find.Execute(
    "aplicacion de python",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "aplicación de python", 2) | Out-Null

# -----------------------------------------------------------------
# 5. Bold + underline the "ACLARACIÓN:" label (splits the run).
# -----------------------------------------------------------------
$aclRange = $d.Content
$aclRange.Find.Execute(
    "ACLARACIÓN:", $true, $false, $false, $false, $false, $true, 1,
    $false, "", 0) | Out-Null
$aclRange.Bold = 1
$aclRange.Underline = 1

Write-Output "done"
